$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assign the brand-new molecule names first (in the order they should be
# interned), regardless of which row they end up in, so the resulting
# shared-string table lists them in this relative order.
$ws.Range("A8").Value = "VANCOMYCIN 1GM"
$ws.Range("A2").Value = "MEROPENEM 1GM"
$ws.Range("A3").Value = "CEFTRIAXONE 1GM"
$ws.Range("A4").Value = "CLINDAMYCIN 600MG"
$ws.Range("A5").Value = "TEICOPLANIN 400MG"
$ws.Range("A6").Value = "LINEZOLID 600MG/300ML"
$ws.Range("A7").Value = "TIGECYCLINE 50MG"
$ws.Range("A9").Value = "DOXYCYCLINE 100MG"

$ws.Range("A6").Select()
